# Unify naming: "Property1" sheet -> "DataNode", "Record" sheet -> "DataTable".
# Also make DataTable (formerly "Record") the active/selected sheet, matching
# the tabSelected/activeTab change in the target workbook.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Property1").Name = "DataNode"
$wb.Worksheets.Item("Record").Name = "DataTable"

$wb.Worksheets.Item("DataTable").Activate()
